# Fruta / hortaliza, semanal
# Insert the latest week's price data (date serial 44628) for
# "Agrícola del Norte S.A. de Arica" / "Betarraga" as two new rows at the
# top of that block (current rows 261-262), pushing the two existing
# blocks (previously rows 261-264) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 261 so the existing
# data (old rows 261-264) shifts down to rows 263-266.
$ws.Rows.Item(261).Resize(2).Insert()

# New "Primera" row (261) - most recent week.
$ws.Cells.Item(261, 1).Value = 1
$ws.Cells.Item(261, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(261, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(261, 4).Value = 44628
$ws.Cells.Item(261, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(261, 5).Value = 15
$ws.Cells.Item(261, 6).Value = 100114014
$ws.Cells.Item(261, 7).Value = "Betarraga"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 700
$ws.Cells.Item(261, 11).Value = 450
$ws.Cells.Item(261, 12).Value = 500
$ws.Cells.Item(261, 13).Value = 475
$ws.Cells.Item(261, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(261, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(261, 16).Value = 119
$ws.Cells.Item(261, 17).Value = 4
$ws.Cells.Item(261, 18).Value = "Hortaliza"

# New "Segunda" row (262) - most recent week.
$ws.Cells.Item(262, 1).Value = 1
$ws.Cells.Item(262, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(262, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(262, 4).Value = 44628
$ws.Cells.Item(262, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(262, 5).Value = 15
$ws.Cells.Item(262, 6).Value = 100114014
$ws.Cells.Item(262, 7).Value = "Betarraga"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Segunda"
$ws.Cells.Item(262, 10).Value = 800
$ws.Cells.Item(262, 11).Value = 450
$ws.Cells.Item(262, 12).Value = 500
$ws.Cells.Item(262, 13).Value = 475
$ws.Cells.Item(262, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(262, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(262, 16).Value = 95
$ws.Cells.Item(262, 17).Value = 5
$ws.Cells.Item(262, 18).Value = "Hortaliza"
